# Auto-generated edit script: update random number cells in sheets x1 and x2
$wb = $excel.ActiveWorkbook

$data1 = @(
    @(2, 0.3612576861019183, -0.4509741943077296),
    @(3, 1.072760912955931, 1.323403035260815),
    @(4, -2.223533206927786, -2.014220626320781),
    @(5, -0.4541552108325183, -1.470128849536629),
    @(6, 0.04363076571425182, 1.31548766476951),
    @(7, 0.5573719852671605, 0.985233437181221),
    @(8, 1.356332413291238, -0.8306242058244614),
    @(9, 1.273320796424101, -0.1623508622203794),
    @(10, 1.921523992568382, -0.1545241888779368),
    @(11, -1.25789201893243, 0.1769685652705896),
    @(12, 0.09143695731503515, -0.762360619456082),
    @(13, 1.263500159971155, 1.857740028936024),
    @(14, -1.577816169294582, 0.6835528390940024),
    @(15, -1.45541334951761, 1.156616098024375),
    @(16, 0.7105966832944799, 1.125894238764907),
    @(17, 2.218016733655032, 0.2757083563123999),
    @(18, -0.6817645475159038, -0.3416254152425416),
    @(19, -0.27493338232481, 1.335834918543516),
    @(20, -0.8486038256819711, -0.6355365040560987),
    @(21, 0.1795837068335649, 0.232123768427802),
    @(22, -0.5910124643022954, -0.5721863496283152),
    @(23, -1.165204771974269, -0.4868852690078531),
    @(24, 0.7544531792817407, 1.160976745851848),
    @(25, -0.1291872541013901, -0.04034992892282423),
    @(26, -0.8425278541924752, -0.9305465810992486),
    @(27, 0.5058555067886076, -1.241718454890224),
    @(28, -0.6100046258041366, 0.4525099422651947),
    @(29, -0.1509779960897119, -0.4423193269248121),
    @(30, -0.7421991892472307, 0.2847661705101655),
    @(31, 0.2804545970828457, 0.9321358610322161),
    @(32, -0.8787503483256203, 0.8894581796848809),
    @(33, -0.05591412642934081, -0.6280738252289281),
    @(34, 1.475409121465419, -3.712209086132272),
    @(35, 1.432567321104183, -1.275052558074382),
    @(36, 0.8313185647595208, -0.6117430478060771),
    @(37, -0.3756477064979983, 0.4507872800818774),
    @(38, 1.262201538123624, 0.500512945526479),
    @(39, -0.07687390305110446, -1.317246099328447),
    @(40, -0.344341875036944, 0.7556265496813992),
    @(41, 0.03240034865184849, 0.7051514399887938),
    @(42, 1.359152849973659, -1.765596695700185),
    @(43, -0.4009145537477865, 0.1840956780260982),
    @(44, -1.630099515024667, 0.6946797792283357),
    @(45, 0.9916997209325101, -0.6492928904908877),
    @(46, -0.3206823915518093, 0.7873452715395616),
    @(47, 0.6391874134074733, -0.9575491450435945),
    @(48, 0.2897337558199167, -0.5943673522018498),
    @(49, -0.05421613005148891, -0.7559319867919181),
    @(50, 0.006279106908708441, -1.082357608942842),
    @(51, -0.9747821915710172, 1.265827095820597),
    @(52, 1.342666537906208, -1.249160017018655),
    @(53, 0.1508419392006338, 1.267987488945343),
    @(54, -1.535235235320922, -0.3855958121579839),
    @(55, 2.415776914015987, 3.640151394197029),
    @(56, 0.09019799498806119, 0.6270576958587428),
    @(57, -0.8540954032864616, 0.2413827670165126),
    @(58, 0.4374416388086468, 0.1006986299289985),
    @(59, -0.9233010094526276, 0.883215672730206),
    @(60, -0.8999821962450502, 1.215670091254053),
    @(61, 0.851753240375213, 0.6351474396423097),
    @(62, 0.7484689341607939, 0.5214220227073247),
    @(63, 0.4935081156732377, 2.128450041800609),
    @(64, -0.03055186249477717, -0.143797569597825),
    @(65, -1.159737680075293, 1.211727073251258),
    @(66, 1.955641840629737, 0.5158039181586543),
    @(67, -0.5615522813127793, -0.7321524595508115),
    @(68, -1.252168502250617, -0.9108252794868594),
    @(69, -0.9926511710365622, 1.144612553394971),
    @(70, -1.073982317310868, -1.088096258958459),
    @(71, -0.1851652111700238, 0.5468494203714616),
    @(72, -1.49915082363527, -0.6895255737874011),
    @(73, -0.6104732474361267, -0.4789128477559996),
    @(74, -0.4851243429616592, -0.9967740784190157),
    @(75, -0.6464803016285395, 0.5059481888541476),
    @(76, -0.2412225223531705, -0.09851401685720333),
    @(77, -0.7589771737214381, 0.612809554523842),
    @(78, -1.504296638388803, 0.4324245218038146),
    @(79, -0.2460990988316427, -0.2305538338288785),
    @(80, 0.7024780289185819, -0.1082448806498865),
    @(81, -0.3701255443416016, 1.614710020551191),
    @(82, 0.6509274161316182, -0.8774883838271481),
    @(83, 1.965929578110975, 0.06149780411508734),
    @(84, -1.124728687083829, -0.5112364039459796),
    @(85, -1.163967608764706, -0.3317493472259974),
    @(86, 1.17192008538866, -1.628158169329231),
    @(87, -0.3409021238902183, -0.8044972596044605),
    @(88, -0.6725688848511649, -1.600808533454595),
    @(89, -2.085134033459898, -0.04795870085326444),
    @(90, 0.1412090957939897, 0.3905611584724428),
    @(91, 0.6909871274527043, 0.7547027560414153),
    @(92, -1.593648091970933, 1.82945583946556),
    @(93, 0.930962100479273, 0.331799236243687),
    @(94, -1.33779385982084, 1.016912381781475),
    @(95, 0.2501472464414053, -0.3899216185163425),
    @(96, -0.813501982848896, -0.806083392596592),
    @(97, 1.755208447096476, -1.829610632654696),
    @(98, -0.7069807519600606, 0.5192920744534519),
    @(99, -1.923254508564905, -1.12854224711177),
    @(100, 0.3361507073238273, -0.7746504858412326),
    @(101, 0.265972173313786, 1.009115080395798)
)

$data2 = @(
    @(2, -0.3105346377330399, 0.05413794446748554),
    @(3, 0.5302369519649347, 1.023604089746236),
    @(4, -0.3147239913714737, -1.695282160566164),
    @(5, -0.6394985134357829, -0.4998963695987991),
    @(6, 0.3788442196684324, -0.9665091379656982),
    @(7, -0.08835473213005463, -0.5831593033059375),
    @(8, 0.1445391582816939, 0.1894315089882786),
    @(9, -0.6106681101487853, 0.93590113117422),
    @(10, -1.082724536445241, 1.846750269793585),
    @(11, 2.393906178139114, 1.048650369095102),
    @(12, 1.333680588244507, 0.1924785980312658),
    @(13, 0.4270434751307727, 1.082913654666637),
    @(14, 0.3318992485700608, -0.01876243590047583),
    @(15, 0.4251144345187799, 0.4879156769374696),
    @(16, 0.4954357930324635, 0.3647759752077492),
    @(17, 0.08025189153805912, 0.7556141739682873),
    @(18, 1.049154377482308, -0.1954352659640649),
    @(19, -0.07861347180143066, -0.7419153130682609),
    @(20, -1.391935024954457, -1.778526639008862),
    @(21, 1.022248024335831, 0.8720867488142698),
    @(22, 0.4844061906485407, -1.003121344778884),
    @(23, -0.1772686227894866, -1.284668345576415),
    @(24, 0.1395888299273537, -0.5270875609555552),
    @(25, 0.1113553175756821, 0.07869917660029876),
    @(26, 1.049698033468303, 0.5431833854864841),
    @(27, 1.302499141139639, 0.2419745310100135),
    @(28, 0.5981377106559861, 0.3378305810398273),
    @(29, 0.271490136353227, -1.429502752229888),
    @(30, -0.6604201317481259, 0.87697278767048),
    @(31, -1.833246096336599, 1.830184171056839),
    @(32, -0.1194446008378872, -0.8296854637853195),
    @(33, -1.403755047936611, 0.6325102522582333),
    @(34, 0.755261757469209, -0.0009800446630148472),
    @(35, 0.05466183204124758, 0.9764544594614774),
    @(36, -1.439140242772982, 0.2658620767454772),
    @(37, -0.6188016485141874, -0.1542412130922351),
    @(38, -0.9065641081416944, 0.09849924825570704),
    @(39, 1.8730419003978, 0.8460741311283735),
    @(40, -0.003917778619547748, 0.4654316113958095),
    @(41, -0.5746793177616006, 0.2460196716699576),
    @(42, -1.281709021113018, -0.2922875756581422),
    @(43, -0.5099091349101297, 0.3551881216036321),
    @(44, 2.156876539855145, 0.01471333943230888),
    @(45, -0.2274977745228388, 0.9751775118699075),
    @(46, 0.9358053671640869, 0.4743102637281847),
    @(47, -0.3447974902068832, -0.5055410197431597),
    @(48, 0.8102855503885027, 0.6976734840309503),
    @(49, 2.732017126720732, 0.4209881714830574),
    @(50, -1.481287607705471, 0.2068318257361759),
    @(51, -0.2373617799508041, 1.093238316293554),
    @(52, -1.54956208507839, -1.007458535687109),
    @(53, 0.04210623463414916, 1.900553891081147),
    @(54, 1.106066249031037, 0.6330632129040724),
    @(55, -0.3755492325763147, -1.724103916176558),
    @(56, -0.6326883428057873, 1.260719207934761),
    @(57, -1.230451127673553, 1.815848526439576),
    @(58, -0.6343662623507739, 0.740857308118788),
    @(59, -0.08900206385252739, -0.4893906338820531),
    @(60, -0.2205669866777333, -0.7254149468221853),
    @(61, 0.7400027493427466, 0.2454393937833157),
    @(62, -0.4352752027077238, 0.8738843205854299),
    @(63, 1.272567391509297, 0.9713454205131548),
    @(64, 1.023093964086024, 0.5923800222623211),
    @(65, 0.9227212165546715, 0.8990875450224115),
    @(66, 3.187003945307631, 0.6589457411383732),
    @(67, 2.933017819675014, 1.16124068144887),
    @(68, 0.008228867749624802, 0.9468593859405309),
    @(69, 1.752863331958491, -0.2132320664740931),
    @(70, 1.285201118173074, -0.3411798483193713),
    @(71, -0.3154736493830692, 0.9866472063630689),
    @(72, -0.6337689394807939, 1.531101448345485),
    @(73, 0.2651550333251318, 1.559658431070031),
    @(74, 0.1872997526973422, -0.2466986571864988),
    @(75, 2.451830903512643, 0.5332391879809889),
    @(76, -0.4756884755098194, -2.354511389150562),
    @(77, 0.154461080228748, -1.640994596724842),
    @(78, -1.828783702619473, 0.3498422223650346),
    @(79, 0.5604316865191199, 0.7358888902569546),
    @(80, -0.6401429203073166, -1.415202937860237),
    @(81, -2.136267171940676, -0.7901796889360349),
    @(82, 1.515104387165458, 0.6682648725172098),
    @(83, -0.4374995114141306, -1.923716795245491),
    @(84, 0.5939927086245168, 0.4102438283422876),
    @(85, 0.448994451571917, 1.01695508869954),
    @(86, 1.633079679120075, -0.5542927483384702),
    @(87, -0.3908248157174969, -0.5233546252828485),
    @(88, 0.2595001950005325, 1.430258011502492),
    @(89, 0.2157158672040455, -0.7478232182979032),
    @(90, 0.1878961303097193, -0.8307693486456922),
    @(91, 0.08954575006630955, -1.608726671336599),
    @(92, -0.6055187808239338, -1.768933128321152),
    @(93, 1.041671334929855, 0.5025018986441411),
    @(94, 0.03220992006580174, -2.044313594589525),
    @(95, -0.4669524915728909, -1.990844805001939),
    @(96, -1.103705625399801, 0.39153116317902),
    @(97, -0.5524511890103677, 0.4071228701456703),
    @(98, 0.1452467142520832, -0.2468653391905227),
    @(99, -0.9586024734149173, 1.952078469304985),
    @(100, 0.6891972536822433, -0.3689784571680829),
    @(101, 0.2303170767222976, 0.3533459810969955)
)

$ws1 = $wb.Worksheets.Item("x1")
foreach ($row in $data1) {
    $r = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
}

$ws2 = $wb.Worksheets.Item("x2")
foreach ($row in $data2) {
    $r = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
}
